# Auto-generated edit script applying the Sephirot_Profits price/profit recompute
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4985.4287
$ws.Range("I19").Value = 5649.8335
$ws.Range("K19").Value = 5649.8335
$ws.Range("M19").Value = -5474.8335
$ws.Range("H98").Value = 1290.75
$ws.Range("I98").Value = 1223.5
$ws.Range("J98").Value = 1492.5
$ws.Range("K98").Value = 1223.5
$ws.Range("L98").Value = 1492.5
$ws.Range("M98").Value = 274.5
$ws.Range("N98").Value = -4488.5
$ws.Range("H118").Value = 474.5
$ws.Range("I118").Value = 474.5
$ws.Range("K118").Value = 1423.5
$ws.Range("M118").Value = 233.5
$ws.Range("H122").Value = 1290.75
$ws.Range("I122").Value = 1223.5
$ws.Range("J122").Value = 1492.5
$ws.Range("K122").Value = 3670.5
$ws.Range("L122").Value = 4477.5
$ws.Range("M122").Value = -1220.5
$ws.Range("N122").Value = -9377.5
$ws.Range("H132").Value = 718.0323
$ws.Range("I132").Value = 708.37933
$ws.Range("K132").Value = 2125.13799
$ws.Range("M132").Value = 404.8620099999998
$ws.Range("H137").Value = 2063.4285
$ws.Range("I137").Value = 1688.8
$ws.Range("K137").Value = 5066.4
$ws.Range("M137").Value = -2516.4
$ws.Range("H138").Value = 2867.8262
$ws.Range("I138").Value = 1269.5555
$ws.Range("J138").Value = 3895.2856
$ws.Range("K138").Value = 3808.6665
$ws.Range("L138").Value = 11685.8568
$ws.Range("M138").Value = 1331.3335
$ws.Range("N138").Value = -21965.8568

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2835.8572
$ws.Range("I61").Value = 1400.25
$ws.Range("J61").Value = 4750
$ws.Range("K61").Value = 1400.25
$ws.Range("L61").Value = 4750
$ws.Range("M61").Value = -1188.25
$ws.Range("N61").Value = -5174
$ws.Range("H95").Value = 29999
$ws.Range("J95").Value = 29999
$ws.Range("L95").Value = 29999
$ws.Range("N95").Value = -35491
$ws.Range("H136").Value = 2835.8572
$ws.Range("I136").Value = 1400.25
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 4200.75
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -1650.75
$ws.Range("N136").Value = -19350

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 64999
$ws.Range("J95").Value = 64999
$ws.Range("L95").Value = 64999
$ws.Range("N95").Value = -70491

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4346.143
$ws.Range("I99").Value = 3606.25
$ws.Range("J99").Value = 5332.6665
$ws.Range("K99").Value = 3606.25
$ws.Range("L99").Value = 5332.6665
$ws.Range("M99").Value = -2108.25
$ws.Range("N99").Value = -8328.666499999999
$ws.Range("H122").Value = 2729.7778
$ws.Range("I122").Value = 2308
$ws.Range("J122").Value = 2940.6667
$ws.Range("K122").Value = 6924
$ws.Range("L122").Value = 8822.000100000001
$ws.Range("M122").Value = -4474
$ws.Range("N122").Value = -13722.0001
$ws.Range("H126").Value = 4346.143
$ws.Range("I126").Value = 3606.25
$ws.Range("J126").Value = 5332.6665
$ws.Range("K126").Value = 10818.75
$ws.Range("L126").Value = 15997.9995
$ws.Range("M126").Value = -8348.75
$ws.Range("N126").Value = -20937.9995
$ws.Range("H132").Value = 3204
$ws.Range("I132").Value = 2845
$ws.Range("K132").Value = 8535
$ws.Range("M132").Value = -6005

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1249.5
$ws.Range("I63").Value = 1249.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3748.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2999.5
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 1218
$ws.Range("I64").Value = 1200
$ws.Range("J64").Value = 1245
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 3735
$ws.Range("M64").Value = -3330
$ws.Range("N64").Value = -4275
$ws.Range("H66").Value = 1249.5
$ws.Range("I66").Value = 1249.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11245.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7501.5
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 1218
$ws.Range("I67").Value = 1200
$ws.Range("J67").Value = 1245
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 3735
$ws.Range("M67").Value = -2664
$ws.Range("N67").Value = -5607
$ws.Range("H81").Value = 2700
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 9000
$ws.Range("N81").Value = -11246
$ws.Range("H84").Value = 2700
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 27000
$ws.Range("N84").Value = -38232
$ws.Range("H121").Value = 993.3333
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 993.3333
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2979.9999
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5599.9999
$ws.Range("H129").Value = 797.1818
$ws.Range("I129").Value = 797.1818
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2391.5454
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2608.4546
$ws.Range("N129").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1957.875
$ws.Range("I102").Value = 1523.2858
$ws.Range("K102").Value = 1523.2858
$ws.Range("M102").Value = 98.71419999999989
$ws.Range("H113").Value = 4374.375
$ws.Range("I113").Value = 4285
$ws.Range("K113").Value = 4285
$ws.Range("M113").Value = -2115
$ws.Range("H122").Value = 8250
$ws.Range("I122").Value = 8250
$ws.Range("K122").Value = 24750
$ws.Range("M122").Value = -22300
$ws.Range("H124").Value = 5555
$ws.Range("J124").Value = 5555
$ws.Range("L124").Value = 5555
$ws.Range("N124").Value = -15375
$ws.Range("H132").Value = 2923.889
$ws.Range("I132").Value = 1422.5
$ws.Range("J132").Value = 4125
$ws.Range("K132").Value = 4267.5
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -1737.5
$ws.Range("N132").Value = -17435

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2408.3333
$ws.Range("I46").Value = 1890
$ws.Range("K46").Value = 1890
$ws.Range("M46").Value = -1702
$ws.Range("H97").Value = 50344
$ws.Range("J97").Value = 50344
$ws.Range("L97").Value = 50344
$ws.Range("N97").Value = -52326
$ws.Range("H132").Value = 3109
$ws.Range("I132").Value = 1617
$ws.Range("J132").Value = 4974
$ws.Range("K132").Value = 4851
$ws.Range("L132").Value = 14922
$ws.Range("M132").Value = -2321
$ws.Range("N132").Value = -19982

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 50020000
$ws.Range("J94").Value = 50020000
$ws.Range("L94").Value = 50020000
$ws.Range("N94").Value = -50021802
$ws.Range("H107").Value = 2357.4285
$ws.Range("I107").Value = 1083.1666
$ws.Range("K107").Value = 3249.4998
$ws.Range("M107").Value = -1329.4998
